$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.044.77'
$ws.Range("E2").Value = '  -4.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.073.24'
$ws.Range("E3").Value = '  -4.93%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '540.78'
$ws.Range("E5").Value = '  -6.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.39'
$ws.Range("E6").Value = '  -11.58%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.068.70'
$ws.Range("E8").Value = '  -4.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -4.62%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.155'
$ws.Range("E10").Value = '  -4.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.21'
$ws.Range("E11").Value = '  -12.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -6.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.40'
$ws.Range("E14").Value = '  -9.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.538.79'
$ws.Range("E15").Value = '  -5.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.977.35'
$ws.Range("E16").Value = '  -4.79%  '

$ws.Range("E17").Value = '  -3.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.079.04'
$ws.Range("E18").Value = '  -5.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.59'
$ws.Range("E19").Value = '  -7.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.49'
$ws.Range("E20").Value = '  -10.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("E21").Value = '  -8.70%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.701'
$ws.Range("E22").Value = '  -5.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.17'
$ws.Range("E23").Value = '  -7.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.29'
$ws.Range("E24").Value = '  -3.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.04'
$ws.Range("E25").Value = '  -10.74%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.33%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  -8.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.15'
$ws.Range("E28").Value = '  -11.75%  '

$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '25.92'
$ws.Range("E30").Value = '  -6.30%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.90'
$ws.Range("E31").Value = '  -15.93%  '

$ws.Range("E32").Value = '  -6.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '58.48'
$ws.Range("E33").Value = '  +6.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.41'
$ws.Range("E34").Value = '  -11.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.96'
$ws.Range("E35").Value = '  -6.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.16'
$ws.Range("E36").Value = '  -7.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '462.57'
$ws.Range("E37").Value = '  -17.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.123.02'
$ws.Range("E38").Value = '  -2.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0389'
$ws.Range("E39").Value = '  -14.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0790'
$ws.Range("E40").Value = '  -8.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.115'
$ws.Range("E41").Value = '  -10.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.03'
$ws.Range("E42").Value = '  -6.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.52'
$ws.Range("E43").Value = '  -12.35%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.249'
$ws.Range("E45").Value = '  -12.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.02'
$ws.Range("E46").Value = '  -13.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.31'
$ws.Range("E47").Value = '  -7.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.77'
$ws.Range("E48").Value = '  -4.76%  '

$ws.Range("E49").Value = '  -4.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0508'
$ws.Range("E50").Value = '  -7.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.99'
$ws.Range("E51").Value = '  -9.34%  '
